# TIYCS_battery BOM update: "new gerber" revision
#  - rename sheet/title  TIYCS_battery -> TIYCS_battery_new
#  - R3 changed from 3.3K/805/C26010  to  1.5K/603/C22843
#  - drop the old 10K R1,R2 resistor row
#  - add a new XT60PW-M connector (J4) row
#  - add a new 13K R2 resistor row
#  - drop the old 1N4148 D3 diode row
#  - BOM table shrinks from G to E (unused F/G columns removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- sheet / title rename ----------------------------------------------
$ws.Name = "Werkblad 1 - TIYCS_battery_new"
$ws.Range("A1").Value = "TIYCS_battery_new"

# ---- remove obsolete rows (bottom-up so row numbers of the other one
#      we still need to delete do not shift) --------------------------
$ws.Rows.Item(15).Delete()   # 1N4148 / D3 / D_SOD-323F / C2128
$ws.Rows.Item(6).Delete()    # 10K / R1,R2 / 805 / C17414

# after the two deletes the sheet looks like:
#  3  1uF              C3,C4,C2  805  C28323
#  4  3.3K              R3        805  C26010     (about to be edited)
#  5  Barrel_Jack_Switch J1       ...  C319134
#  6  B540              D2        D_SMC C57579
#  7  220uF/50V         C1        CP_Elec_10x10  C190286
#  8  D                 D1        D_SMA C22452
#  9  External Jack     J2        ...   C8269
# 10  XL4015            U3        CONV_XL4015 C51661
# 11  47uH/5A           L1        YSPI1365  C497913
# 12  330uF/25V         C5        C_Elec_10x10.2 C2687701
# 13  killSwitch        S1        MR1-110-C5N-BB C268222

# ---- insert the two new rows (again in an order that keeps the not-yet
#      -used insertion point stable) ------------------------------------
$ws.Rows.Item(9).Insert()    # makes room for the 13K / R2 row, before "External Jack"
$ws.Rows.Item(8).Insert()    # makes room for the XT60PW-M / J4 row, before "D"

# final row layout is now 3..15, matching the target sheet exactly.

# ---- R3 comment/footprint/LCSC update (row 4) --------------------------
$ws.Cells.Item(4,1).Value = "1.5K"
$ws.Cells.Item(4,3).Value = 603
$ws.Cells.Item(4,4).Value = "C22843"

# ---- new row 8: XT60PW-M connector -------------------------------------
$ws.Cells.Item(8,1).Value = "XT60PW-M"
$ws.Cells.Item(8,2).Value = "J4"
$ws.Cells.Item(8,3).Value = "AMASS_XT60PW-M"
$ws.Cells.Item(8,4).Value = "C98732"
$ws.Cells.Item(8,5).Value = 1

# ---- new row 10: 13K resistor ------------------------------------------
$ws.Cells.Item(10,1).Value = "13K"
$ws.Cells.Item(10,2).Value = "R2"
$ws.Cells.Item(10,3).Value = 805
$ws.Cells.Item(10,4).Value = "C17455"
$ws.Cells.Item(10,5).Value = 1

# ---- remove now-unused trailing F:G columns ----------------------------
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(6).Delete()

# ---- column widths (A/B got wider once F/G disappeared) ----------------
$ws.Columns.Item(1).ColumnWidth = 19.2857142857143
$ws.Columns.Item(2).ColumnWidth = 15.8571428571429

# ---- row heights --------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 27.65
$ws.Rows.Item(2).RowHeight = 14.7
for ($r = 3; $r -le 15; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.7
}
